$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("setup")

# Update the data.path column (C) for rows 2-5 from "../data" to "../example_data"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "../example_data"
}

# Update the selection shown on the setup sheet
$ws.Activate()
$ws.Range("C3:C5").Select()
